# The workbook previously used the blank/base RDF prefix (rendered as ":xxx")
# for the "ome_core" namespace. This change gives that base prefix an explicit
# name, "ome", so every cell that referenced the bare ":xxx" form now reads
# "ome:xxx", and the "@prefix" table gets a new row-1 entry naming the base
# prefix "ome" (column A) alongside its namespace URI (column B, unchanged).

$wb = $excel.ActiveWorkbook

# --- Sheet "@prefix": add the new base-prefix name in A1 ---
$wsPrefix = $wb.Worksheets.Item("@prefix")
$wsPrefix.Range("A1").Value = "ome"

# --- Sheet "Image": rename the ":xxx" references to "ome:xxx" ---
$wsImage = $wb.Worksheets.Item("Image")
$wsImage.Range("E3").Value = "ome:pixels"
$wsImage.Range("F3").Value = "ome:acquisitionDate"
$wsImage.Range("G3").Value = "ome:annotation"
$wsImage.Range("B4").Value = "ome:Image"
$wsImage.Range("E4").Value = "ome:Pixels"
$wsImage.Range("G4").Value = "ome:Annotation"

# --- Sheet "Pixels": rename the ":xxx" references to "ome:xxx" ---
$wsPixels = $wb.Worksheets.Item("Pixels")
$wsPixels.Range("D3").Value = "ome:pixelType"
$wsPixels.Range("E3").Value = "ome:dimensionOrder"
$wsPixels.Range("F3").Value = "ome:physicalSizeX"
$wsPixels.Range("G3").Value = "ome:physicalSizeY"
$wsPixels.Range("H3").Value = "ome:sizeC"
$wsPixels.Range("I3").Value = "ome:sizeT"
$wsPixels.Range("J3").Value = "ome:sizeX"
$wsPixels.Range("K3").Value = "ome:sizeY"
$wsPixels.Range("L3").Value = "ome:sizeZ"
$wsPixels.Range("M3").Value = "ome:binData"
$wsPixels.Range("B4").Value = "ome:Pixels"
$wsPixels.Range("D4").Value = "ome:PixelType"
$wsPixels.Range("E4").Value = "ome:DimensionOrder"
$wsPixels.Range("M4").Value = "ome:BinData"

# --- Sheet "Binary_Data": rename the ":xxx" references to "ome:xxx" ---
$wsBinData = $wb.Worksheets.Item("Binary_Data")
$wsBinData.Range("C3").Value = "ome:bigEndian"
$wsBinData.Range("D3").Value = "ome:data"
$wsBinData.Range("E3").Value = "ome:length"
$wsBinData.Range("B4").Value = "ome:BinData"

# --- Sheet "Structured_Annotations": rename the ":xxx" references to "ome:xxx" ---
$wsStructAnn = $wb.Worksheets.Item("Structured_Annotations")
$wsStructAnn.Range("C3").Value = "ome:annotation"
$wsStructAnn.Range("B4").Value = "ome:StructuredAnnotations"
$wsStructAnn.Range("C4").Value = "ome:Annotation"

# --- Sheet "XML_Annotation": rename the ":xxx" references to "ome:xxx" ---
$wsXmlAnn = $wb.Worksheets.Item("XML_Annotation")
$wsXmlAnn.Range("E3").Value = "ome:nameSpace"
$wsXmlAnn.Range("B4").Value = "ome:XMLAnnotation"

$wb.Save()
